$d = $word.ActiveDocument

$replacements = @(
    @{old = "272÷7="; new = "188÷9="},
    @{old = "800÷6="; new = "266÷4="},
    @{old = "193÷6="; new = "760÷6="},
    @{old = "945÷5="; new = "972÷6="},
    @{old = "554÷7="; new = "282÷5="},
    @{old = "386÷3="; new = "990÷9="},
    @{old = "639÷9="; new = "957÷6="},
    @{old = "199÷2="; new = "697÷4="},
    @{old = "915÷5="; new = "900÷2="},
    @{old = "723÷9="; new = "712÷2="},
    @{old = "541÷2="; new = "833÷4="},
    @{old = "328÷8="; new = "689÷3="},
    @{old = "572÷2="; new = "445÷4="},
    @{old = "467÷2="; new = "373÷7="},
    @{old = "878÷6="; new = "719÷8="},
    @{old = "689÷6="; new = "431÷4="},
    @{old = "247÷6="; new = "144÷2="},
    @{old = "100÷9="; new = "263÷5="},
    @{old = "575÷9="; new = "695÷3="},
    @{old = "903÷8="; new = "931÷2="},
    @{old = "594÷5="; new = "848÷8="},
    @{old = "822÷3="; new = "510÷4="},
    @{old = "794÷8="; new = "958÷6="},
    @{old = "552÷8="; new = "782÷9="},
    @{old = "437÷9="; new = "313÷3="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
